# Apply the "auditee_uei" column insertion + findings-row update to every
# worksheet in the workbook (both "93" and "OVER-93" sheets carry the same
# table/shape in this file).

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- 1. Insert a new column at D ("auditee_uei"), shifting
    #         award_reference..prior_finding_ref_numbers one column right
    #         (D:U -> E:V). Excel carries formatting/column-widths along
    #         automatically for the existing columns.
    $ws.Columns.Item(4).Insert()

    # New header + value for the inserted auditee_uei column.
    $ws.Range("D1").Value = "auditee_uei"
    $ws.Range("D2").Value = "ZHXHKTS56XM1"

    # Column D / new column I (federal_program_name, now shorter text)
    # get narrower custom widths (closest achievable via the COM
    # pixel-rounded ColumnWidth property to the authored 16.8 / 26.4).
    $ws.Range("D1").ColumnWidth = 16
    $ws.Range("I1").ColumnWidth = 25.5

    # --- 2. Update the finding row's data values that changed.

    # aln (now column G): "93.778" -> "93.498", staying a plain text
    # value instead of being auto-converted to a number. Enter it as a
    # string-literal formula, then collapse it back to a plain value so
    # no formula (and no extra number-format style) is left behind.
    $ws.Range("G2").Formula = "=""93.498"""
    $ws.Range("G2").Copy()
    $ws.Range("G2").PasteSpecial(-4163)
    $ws.Application.CutCopyMode = $false

    # federal_program_name (now column I)
    $ws.Range("I2").Value = "PROVIDER RELIEF FUND"

    # amount_expended (now column J)
    $ws.Range("J2").Value = 2485265

    # is_direct (now column K): NO -> YES, highlighted like the
    # existing is_other_matters / is_material_weakness cells (copy the
    # fill from the already-highlighted P2 cell, which reuses the
    # existing style instead of creating a duplicate one).
    $ws.Range("K2").Value = "YES"
    $ws.Range("P2").Copy()
    $ws.Range("K2").PasteSpecial(-4122)

    # is_major (now column L): NO -> YES, highlighted.
    $ws.Range("L2").Value = "YES"
    $ws.Range("P2").Copy()
    $ws.Range("L2").PasteSpecial(-4122)

    $ws.Application.CutCopyMode = $false

    # is_passthrough_award (now column M): "NO" text -> boolean FALSE.
    $ws.Range("M2").Value = $false

    # is_repeat_finding (now column U): "NO" text -> boolean FALSE.
    $ws.Range("U2").Value = $false
}
